$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D. This shifts the existing D:K data
# (and their formatting) one column to the right, into E:L, and extends
# the sheet from A5:K102 to A5:L102.
$ws.Columns("D").Insert()

# The newly inserted column D has no formatting yet. Copy the number
# formats/styles from the (now shifted) column E, which has the same
# per-row style pattern (date header row uses style 2, data rows use
# style 3, blank separator rows use style 3 with no value, etc.)
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the new (most recent) period's figures.

$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 297900
$ws.Range("D9").Value = 89700
$ws.Range("D10").Value = 208300
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 7500
$ws.Range("D15").Value = 114300
$ws.Range("D17").Value = 247600
$ws.Range("D18").Value = 50400
$ws.Range("D20").Value = 7700
$ws.Range("D21").Value = 169200
$ws.Range("D22").Value = 56500
$ws.Range("D23").Value = 1500
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 1500
$ws.Range("D27").Value = 9200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -7700
$ws.Range("D33").Value = 9200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 9200
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 17100
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 90400
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 20900
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 1071400
$ws.Range("D48").Value = 2860400
$ws.Range("D49").Value = 34100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 77400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 4174000
$ws.Range("D57").Value = 49100
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 23300
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 1323800
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1433000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -1684100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2741000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 9200
$ws.Range("D83").Value = 111200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 97600
$ws.Range("D91").Value = -71000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 375700
$ws.Range("D96").Value = -95100
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -465800
$ws.Range("D101").Value = "NA"
$ws.Range("D102").Value = 7600

# A few rows were also restated for the next three (shifted) periods,
# not only the brand-new column - update those explicitly.
$ws.Range("E20").Value = 9500
$ws.Range("F20").Value = 7300
$ws.Range("G20").Value = 6100
$ws.Range("E21").Value = 170500
$ws.Range("F21").Value = 248600
$ws.Range("G21").Value = 317100
$ws.Range("E23").Value = -2300
$ws.Range("F23").Value = 20000
$ws.Range("G23").Value = 22300
$ws.Range("E26").Value = -2100
$ws.Range("F26").Value = 19500
$ws.Range("G26").Value = 21900
$ws.Range("E32").Value = -9500
$ws.Range("F32").Value = -7300
$ws.Range("G32").Value = -6100
